$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Monthly date headers (Jan-2019 .. Dec-2022), as Excel serial date numbers
$dates = @(43466,43497,43525,43556,43586,43617,43647,43678,43709,43739,43770,43800,43831,43862,43891,43922,43952,43983,44013,44044,44075,44105,44136,44166,44197,44228,44256,44287,44317,44348,44378,44409,44440,44470,44501,44531,44562,44593,44621,44652,44682,44713,44743,44774,44805,44835,44866,44896)

# Replace the text month headers in row 1 and row 21 with real dates
# formatted as mmm-yy, wiping any pre-existing per-cell formatting first
# so every header cell lands on a single uniform style.
$ws.Range("B1:AW1").ClearFormats()
$ws.Range("B21:AW21").ClearFormats()

for ($i = 0; $i -lt $dates.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $dates[$i]
    $ws.Cells.Item(21, $col).Value = $dates[$i]
}

$ws.Range("B1:AW1").NumberFormat = "mmm-yy"
$ws.Range("B21:AW21").NumberFormat = "mmm-yy"

# Re-enter the Grand Total sums across the row in one go so Excel
# collapses them back into a single shared formula, same as the original.
$ws.Range("B19:AW19").Formula = "=SUM(B2:B18)"

# Leave the cursor where the author last clicked before saving.
$ws.Range("I35").Select()
